$d = $word.ActiveDocument

# Locate the paragraph that holds the "m:...asImage().fit(...)" field
# (the one built from w:fldChar / w:instrText runs).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $target = $p
    }
}

if ($target -eq $null) {
    $target = $d.Paragraphs(2)
}

$start = $target.Range.Start
$end = $target.Range.End - 1

$r = $d.Range($start, $end)

# Rebuild the paragraph: the field-code runs (w:fldChar begin/end +
# w:instrText) become plain w:t runs, the field delimiter runs become
# literal "{" / "}" text runs, and the two space-only instrText runs
# (right after the opening brace and right before the closing one) are
# dropped, matching the TokenIteratorFieldRewriterSplit output.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="2F8A187F" w14:textId="0CC007EB" w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">
<w:r><w:t>{</w:t></w:r>
<w:r w:rsidR="00DE6D5A"><w:t>m</w:t></w:r>
<w:r><w:t>:</w:t></w:r>
<w:r w:rsidR="004B598D"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>'</w:t></w:r>
<w:r w:rsidR="00D67687" w:rsidRPr="00D67687"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>Mona_Lisa.jpg</w:t></w:r>
<w:r w:rsidR="004B598D"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>'.asImage()</w:t></w:r>
<w:r w:rsidR="0047710F"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>.fit(</w:t></w:r>
<w:r w:rsidR="00355CDE"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>1</w:t></w:r>
<w:r w:rsidR="007C11BF"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>0</w:t></w:r>
<w:r w:rsidR="0047710F"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t xml:space="preserve">0, </w:t></w:r>
<w:r w:rsidR="007C11BF"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>15</w:t></w:r>
<w:r w:rsidR="0047710F"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>0</w:t></w:r>
<w:r w:rsidR="00C9704C"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>, false</w:t></w:r>
<w:r w:rsidR="0047710F"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>)</w:t></w:r>
<w:r><w:t xml:space="preserve">}</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r.InsertXML($xml)
